$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.593.22"
$ws.Range("E2").Value = "  -8.59%  "
$ws.Range("D3").Value = "'1.640.65"
$ws.Range("E3").Value = "  -9.94%  "
$ws.Range("D4").Value = "'1.018"
$ws.Range("E4").Value = "  +1.20%  "
$ws.Range("D5").Value = "'218.77"
$ws.Range("E5").Value = "  -6.16%  "
$ws.Range("D6").Value = "'1.018"
$ws.Range("E6").Value = "  +1.25%  "
$ws.Range("D7").Value = "'0.5015"
$ws.Range("E7").Value = "  -15.00%  "
$ws.Range("D8").Value = "'0.2516"
$ws.Range("E8").Value = "  -8.21%  "
$ws.Range("D9").Value = "'21.49"
$ws.Range("E9").Value = "  -6.59%  "
$ws.Range("D10").Value = "'0.06063"
$ws.Range("E10").Value = "  -10.79%  "
$ws.Range("D11").Value = "'0.07359"
$ws.Range("E11").Value = "  -1.96%  "
$ws.Range("D12").Value = "'1.653.35"
$ws.Range("E12").Value = "  -13.47%  "
$ws.Range("D13").Value = "'4.457"
$ws.Range("E13").Value = "  -4.42%  "
$ws.Range("D14").Value = "'0.5676"
$ws.Range("E14").Value = "  -8.90%  "
$ws.Range("D15").Value = "'1.869.19"
$ws.Range("E15").Value = "  -9.59%  "
$ws.Range("D16").Value = "'0.000007939"
$ws.Range("E16").Value = "  -14.69%  "
$ws.Range("D17").Value = "'63.39"
$ws.Range("E17").Value = "  -14.75%  "
$ws.Range("D18").Value = "'26.608.44"
$ws.Range("E18").Value = "  -7.65%  "
$ws.Range("D19").Value = "'4.933"
$ws.Range("E19").Value = "  -8.85%  "
$ws.Range("D20").Value = "'1.018"
$ws.Range("E20").Value = "  +1.25%  "
$ws.Range("D21").Value = "'10.58"
$ws.Range("E21").Value = "  -6.85%  "
$ws.Range("D22").Value = "'182.55"
$ws.Range("E22").Value = "  -12.02%  "
$ws.Range("D23").Value = "'1.018"
$ws.Range("E23").Value = "  +1.12%  "
$ws.Range("D24").Value = "'6.146"
$ws.Range("E24").Value = "  -9.08%  "
$ws.Range("D25").Value = "'142.42"
$ws.Range("E25").Value = "  -7.46%  "
$ws.Range("D26").Value = "'7.536"
$ws.Range("E26").Value = "  -2.94%  "
$ws.Range("D27").Value = "'0.1138"
$ws.Range("E27").Value = "  -10.11%  "
$ws.Range("D28").Value = "'14.95"
$ws.Range("E28").Value = "  -7.93%  "
$ws.Range("D29").Value = "'1.331"
$ws.Range("E29").Value = "  -5.34%  "
$ws.Range("D30").Value = "'0.05687"
$ws.Range("E30").Value = "  -10.95%  "
$ws.Range("D31").Value = "'1.335"
$ws.Range("E31").Value = "  -6.64%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'3.398"
$ws.Range("E32").Value = "  -8.25%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "'3.401"
$ws.Range("E33").Value = "  -7.26%  "
$ws.Range("D34").Value = "'1.543"
$ws.Range("E34").Value = "  -7.97%  "
$ws.Range("D35").Value = "'0.9663"
$ws.Range("E35").Value = "  -7.86%  "
$ws.Range("D36").Value = "'2.442"
$ws.Range("E36").Value = "  -3.65%  "
$ws.Range("D37").Value = "'0.5886"
$ws.Range("E37").Value = "  -6.78%  "
$ws.Range("D38").Value = "'2.602"
$ws.Range("E38").Value = "  -5.48%  "
$ws.Range("D39").Value = "'0.01562"
$ws.Range("E39").Value = "  -8.11%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'0.8596"
$ws.Range("E40").Value = "  -1.21%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value = "'1.018"
$ws.Range("E41").Value = "  +1.27%  "
$ws.Range("D42").Value = "'1.054.37"
$ws.Range("E42").Value = "  -6.90%  "
$ws.Range("D43").Value = "'5.666"
$ws.Range("E43").Value = "  -12.09%  "
$ws.Range("D44").Value = "'95.02"
$ws.Range("E44").Value = "  -4.53%  "
$ws.Range("D45").Value = "'1.780.16"
$ws.Range("E45").Value = "  -9.88%  "
$ws.Range("D46").Value = "'0.00000000108"
$ws.Range("E46").Value = "  -3.83%  "
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").Value = "'1.009"
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "'0.4409"
$ws.Range("E48").Value = "  -2.63%  "
$ws.Range("D49").Value = "'54.55"
$ws.Range("E49").Value = "  -9.10%  "
$ws.Range("D50").Value = "'0.05218"
$ws.Range("E50").Value = "  -4.79%  "
$ws.Range("D51").Value = "'7.701"
$ws.Range("E51").Value = "  -6.37%  "
